$d = $word.ActiveDocument

# --- Locate the table cell that currently holds "{{pin_code}} {{qr_code}}".
#     (2nd table, row 1, col 2 in before.docx -> document paragraph #99.) ---
$pinPara = $d.Paragraphs.Item(99)
if ($pinPara.Range.Text -notmatch "pin_code") {
    throw "Expected paragraph 99 to contain the pin_code/qr_code placeholders, found: $($pinPara.Range.Text)"
}

# --- Step 1: rebuild that paragraph as "{{pin_code_with_qr: right}}",
#     matching the run/proofErr layout used elsewhere in the template. ---
$xml = '<w:p><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t>pin_code_with_qr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:lang w:val="en-US"/></w:rPr><w:t>right}}</w:t></w:r></w:p>'
[void]$pinPara.Range.InsertXML($xml)

# --- Step 2: remove the now-redundant empty paragraph that used to follow
#     the pin/qr paragraph inside the same table cell, merging the cell
#     back down to a single paragraph. ---
$emptyCellPara = $d.Paragraphs.Item(100)
# A cell's final paragraph reports its range text with a trailing cell-mark
# (chr 7) after the paragraph mark (chr 13), so trim both before checking.
$emptyCellText = $emptyCellPara.Range.Text.TrimEnd([char]13, [char]7)
if ($emptyCellText -ne "") {
    throw "Expected paragraph 100 to be empty, found: $emptyCellText"
}
$emptyCellPara.Range.Delete()

# --- Step 3: the two trailing centered empty paragraphs (right after the
#     table, before the section break) switch language ru-RU -> en-US.
#     Their indices shift down by one now that paragraph 100 is gone. ---
$tail1 = $d.Paragraphs.Item(101)
$tail2 = $d.Paragraphs.Item(102)
if ($tail1.Range.LanguageID -ne "ru-RU" -or $tail2.Range.LanguageID -ne "ru-RU") {
    throw "Expected the two trailing paragraphs to be ru-RU, found: $($tail1.Range.LanguageID) / $($tail2.Range.LanguageID)"
}
$tail1.Range.LanguageID = "en-US"
$tail2.Range.LanguageID = "en-US"
